$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("AlfalahIslamicSpecialRateTDR")
$ws2 = $wb.Worksheets.Item("IslamicSpRateTDRMonthly")
$ws3 = $wb.Worksheets.Item("DealslipforspecialTDR")

# Sheet1 updates
$ws1.Range("C2").Value = "TDR1M"
$ws1.Range("D2").Value = "1M"

# Sheet2 updates
$ws2.Range("A2").Value = 11871224
$ws2.Range("F1").Value = "DRAWDOWN.ACCOUNT"
$ws2.Range("G1").Value = "PRIN.LIQ.ACCT"
$ws2.Range("H1").Value = "INT.LIQ.ACCT"
$ws2.Range("F2:H2").NumberFormat = "@"
$ws2.Range("F2").Value = "5000000539"
$ws2.Range("G2").Value = "5000000539"
$ws2.Range("H2").Value = "5000000539"

# Sheet3 updates
$ws3.Range("A2").Value = 11871234
$ws3.Range("F1").Value = "DRAWDOWN.ACCOUNT"
$ws3.Range("G1").Value = "PRIN.LIQ.ACCT"
$ws3.Range("H1").Value = "INT.LIQ.ACCT"
$ws3.Range("F2:H2").NumberFormat = "@"
$ws3.Range("F2").Value = "5000000540"
$ws3.Range("G2").Value = "5000000540"
$ws3.Range("H2").Value = "5000000540"

# New columns on sheet3 need explicit (best-fit) widths like the rest of the sheet
$ws3.Columns.Item(5).ColumnWidth = 15
$ws3.Columns.Item(6).ColumnWidth = 21.5
$ws3.Columns.Item(7).ColumnWidth = 13.166666666666666
$ws3.Columns.Item(8).ColumnWidth = 11.666666666666666

# Activate sheet1 tab
$ws1.Activate()

# Selections
$ws1.Range("D17").Select()
$ws2.Range("A2").Select()
$ws3.Range("E21").Select()
$ws1.Activate()
